# Auto-generated Excel COM-interop script
# Applies scheduled market-price / profit recompute updates across all 8 sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7372.909
$ws.Range("J40").Value = 9999
$ws.Range("L40").Value = 9999
$ws.Range("N40").Value = -10349
$ws.Range("H76").Value = 3003
$ws.Range("I76").Value = 3003
$ws.Range("K76").Value = 3003
$ws.Range("M76").Value = -2688
$ws.Range("H79").Value = 3003
$ws.Range("I79").Value = 3003
$ws.Range("K79").Value = 3003
$ws.Range("M79").Value = -1911
$ws.Range("H88").Value = 1179.6666
$ws.Range("J88").Value = 1179.6666
$ws.Range("L88").Value = 1179.6666
$ws.Range("N88").Value = -1991.6666
$ws.Range("H91").Value = 1179.6666
$ws.Range("J91").Value = 1179.6666
$ws.Range("L91").Value = 1179.6666
$ws.Range("N91").Value = -3987.6666
$ws.Range("H112").Value = 2467.5
$ws.Range("J112").Value = 1826.25
$ws.Range("L112").Value = 5478.75
$ws.Range("N112").Value = -7694.75
$ws.Range("H113").Value = 7596.143
$ws.Range("I113").Value = 7374.6665
$ws.Range("K113").Value = 7374.6665
$ws.Range("M113").Value = -4120.6665
$ws.Range("H132").Value = 2310
$ws.Range("I132").Value = 1307.3077
$ws.Range("K132").Value = 3921.9231
$ws.Range("M132").Value = -1391.9231
$ws.Range("H137").Value = 1829.9166
$ws.Range("I137").Value = 915.7143
$ws.Range("K137").Value = 2747.1429
$ws.Range("M137").Value = -197.1428999999998
$ws.Range("H138").Value = 2589.875
$ws.Range("I138").Value = 1580.6364
$ws.Range("K138").Value = 4741.9092
$ws.Range("M138").Value = 398.0907999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1047.238
$ws.Range("I32").Value = 946.375
$ws.Range("J32").Value = 1370
$ws.Range("K32").Value = 946.375
$ws.Range("L32").Value = 1370
$ws.Range("M32").Value = -659.375
$ws.Range("N32").Value = -1944
$ws.Range("H61").Value = 1728.5
$ws.Range("I61").Value = 1728.5
$ws.Range("K61").Value = 1728.5
$ws.Range("M61").Value = -1516.5
$ws.Range("H63").Value = 8797.799999999999
$ws.Range("J63").Value = 10247.25
$ws.Range("L63").Value = 10247.25
$ws.Range("N63").Value = -11619.25
$ws.Range("H66").Value = 8797.799999999999
$ws.Range("J66").Value = 10247.25
$ws.Range("L66").Value = 51236.25
$ws.Range("N66").Value = -58100.25
$ws.Range("H74").Value = 3498.889
$ws.Range("I74").Value = 3281.4167
$ws.Range("K74").Value = 3281.4167
$ws.Range("M74").Value = -2407.4167
$ws.Range("H77").Value = 3498.889
$ws.Range("I77").Value = 3281.4167
$ws.Range("K77").Value = 16407.0835
$ws.Range("M77").Value = -12039.0835
$ws.Range("H88").Value = 2012
$ws.Range("I88").Value = 1510
$ws.Range("K88").Value = 1510
$ws.Range("M88").Value = -1104
$ws.Range("H91").Value = 2012
$ws.Range("I91").Value = 1510
$ws.Range("K91").Value = 1510
$ws.Range("M91").Value = -106
$ws.Range("H136").Value = 1728.5
$ws.Range("I136").Value = 1728.5
$ws.Range("K136").Value = 5185.5
$ws.Range("M136").Value = -2635.5
$ws.Range("H138").Value = 99992.5
$ws.Range("J138").Value = 99992.5
$ws.Range("L138").Value = 99992.5
$ws.Range("N138").Value = -110272.5
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3006.375
$ws.Range("I86").Value = 1450.25
$ws.Range("K86").Value = 1450.25
$ws.Range("M86").Value = -327.25
$ws.Range("H89").Value = 3006.375
$ws.Range("I89").Value = 1450.25
$ws.Range("K89").Value = 7251.25
$ws.Range("M89").Value = -1635.25
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 74.40000000000001
$ws.Range("I7").Value = 74.59999999999999
$ws.Range("K7").Value = 74.59999999999999
$ws.Range("M7").Value = 38.40000000000001
$ws.Range("H22").Value = 870.7646999999999
$ws.Range("I22").Value = 779.4167
$ws.Range("J22").Value = 1090
$ws.Range("K22").Value = 779.4167
$ws.Range("L22").Value = 1090
$ws.Range("M22").Value = -429.4167
$ws.Range("N22").Value = -1790
$ws.Range("H51").Value = 51994
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 51994
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 51994
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -53466
$ws.Range("H61").Value = 51994
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 51994
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 51994
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -52690
$ws.Range("H99").Value = 4105.4
$ws.Range("J99").Value = 3842.3333
$ws.Range("L99").Value = 3842.3333
$ws.Range("N99").Value = -6838.3333
$ws.Range("H126").Value = 4105.4
$ws.Range("J126").Value = 3842.3333
$ws.Range("L126").Value = 11526.9999
$ws.Range("N126").Value = -16466.9999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 333908.12
$ws.Range("I4").Value = 375646.25
$ws.Range("K4").Value = 1126938.75
$ws.Range("M4").Value = -1126826.75
$ws.Range("H26").Value = 151.16667
$ws.Range("I26").Value = 151.16667
$ws.Range("K26").Value = 453.50001
$ws.Range("M26").Value = -165.50001
$ws.Range("H44").Value = 5223.8335
$ws.Range("I44").Value = 147.66667
$ws.Range("J44").Value = 10300
$ws.Range("K44").Value = 443.00001
$ws.Range("L44").Value = 30900
$ws.Range("M44").Value = -45.00001000000003
$ws.Range("N44").Value = -31696
$ws.Range("H58").Value = 5695
$ws.Range("I58").Value = 5695
$ws.Range("K58").Value = 17085
$ws.Range("M58").Value = -16957
$ws.Range("H140").Value = 3384.3635
$ws.Range("I140").Value = 2824.375
$ws.Range("K140").Value = 8473.125
$ws.Range("M140").Value = -3293.125

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 432.35715
$ws.Range("I2").Value = 172.55556
$ws.Range("K2").Value = 172.55556
$ws.Range("M2").Value = -59.55556000000001
$ws.Range("H113").Value = 8000
$ws.Range("I113").Value = 8000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 8000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -5830
$ws.Range("N113").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1005000
$ws.Range("I2").Value = 1002500
$ws.Range("K2").Value = 1002500
$ws.Range("M2").Value = -1002388
$ws.Range("H22").Value = 692.8333
$ws.Range("I22").Value = 225.5
$ws.Range("K22").Value = 225.5
$ws.Range("M22").Value = 69.5
$ws.Range("H27").Value = 692.8333
$ws.Range("I27").Value = 225.5
$ws.Range("K27").Value = 225.5
$ws.Range("M27").Value = -118.5
$ws.Range("H48").Value = 5036
$ws.Range("I48").Value = 5036
$ws.Range("K48").Value = 5036
$ws.Range("M48").Value = -4375
$ws.Range("H55").Value = 1444.8334
$ws.Range("I55").Value = 1783
$ws.Range("J55").Value = 1203.2858
$ws.Range("K55").Value = 1783
$ws.Range("L55").Value = 1203.2858
$ws.Range("M55").Value = -1610
$ws.Range("N55").Value = -1549.2858
$ws.Range("H82").Value = 3435.5386
$ws.Range("I82").Value = 568.8333
$ws.Range("K82").Value = 568.8333
$ws.Range("M82").Value = -207.8333
$ws.Range("H85").Value = 3435.5386
$ws.Range("I85").Value = 568.8333
$ws.Range("K85").Value = 568.8333
$ws.Range("M85").Value = 679.1667
$ws.Range("H93").Value = 1306.9375
$ws.Range("I93").Value = 1387.2858
$ws.Range("J93").Value = 1244.4445
$ws.Range("K93").Value = 1387.2858
$ws.Range("L93").Value = 1244.4445
$ws.Range("M93").Value = -139.2858000000001
$ws.Range("N93").Value = -3740.4445
$ws.Range("H122").Value = 7038.2
$ws.Range("I122").Value = 7038.2
$ws.Range("K122").Value = 21114.6
$ws.Range("M122").Value = -18664.6
$ws.Range("H132").Value = 4928.143
$ws.Range("I132").Value = 4928.143
$ws.Range("K132").Value = 14784.429
$ws.Range("M132").Value = -12254.429

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1559.2667
$ws.Range("I2").Value = 1755.3846
$ws.Range("K2").Value = 1755.3846
$ws.Range("M2").Value = -1643.3846
$ws.Range("H115").Value = 40000
$ws.Range("I115").Value = 20000
$ws.Range("J115").Value = 60000
$ws.Range("K115").Value = 20000
$ws.Range("L115").Value = 60000
$ws.Range("M115").Value = -18433
$ws.Range("N115").Value = -63134
$ws.Range("H136").Value = 3528.9
$ws.Range("I136").Value = 2499.5715
$ws.Range("K136").Value = 7498.7145
$ws.Range("M136").Value = -4948.7145
